# Natmi following Dr Hou advice:
# Recompute the Cthrc1(ligand, sender)->Ror2(receptor, target) edge table for
# all 3x3 sender/target cluster combinations (ECs, FAPs, sCs), replacing the
# previous 3-row table (one row per sender, fixed target) with the full 9-row
# cross product and updated statistics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cthrc1"
$ws.Range("C2").Value = "Ror2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.507882
$ws.Range("H2").Value = 16.523646
$ws.Range("I2").Value = 0.03518866199235487
$ws.Range("J2").Value = 0.03518866199235487
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1378523333333333
$ws.Range("N2").Value = 0.413557
$ws.Range("O2").Value = 0.01719056794796269
$ws.Range("P2").Value = 0.01719056794796269
$ws.Range("Q2").Value = 0.7592743854246665
$ws.Range("R2").Value = 6.833469468822
$ws.Range("S2").Value = 0.0006049130849774687
$ws.Range("T2").Value = 0.0006049130849774687
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cthrc1"
$ws.Range("C3").Value = "Ror2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.507882
$ws.Range("H3").Value = 16.523646
$ws.Range("I3").Value = 0.03518866199235487
$ws.Range("J3").Value = 0.03518866199235487
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 7.692787333333334
$ws.Range("N3").Value = 23.078362
$ws.Range("O3").Value = 0.9593118967607375
$ws.Range("P3").Value = 0.9593118967607375
$ws.Range("Q3").Value = 42.37096488309466
$ws.Range("R3").Value = 381.338683947852
$ws.Range("S3").Value = 0.03375690208035842
$ws.Range("T3").Value = 0.03375690208035842
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Cthrc1"
$ws.Range("C4").Value = "Ror2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5.507882
$ws.Range("H4").Value = 16.523646
$ws.Range("I4").Value = 0.03518866199235487
$ws.Range("J4").Value = 0.03518866199235487
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1884283333333333
$ws.Range("N4").Value = 0.565285
$ws.Range("O4").Value = 0.02349753529129985
$ws.Range("P4").Value = 0.02349753529129985
$ws.Range("Q4").Value = 1.037841025456667
$ws.Range("R4").Value = 9.340569229110001
$ws.Range("S4").Value = 0.0008268468270189803
$ws.Range("T4").Value = 0.0008268468270189803
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cthrc1"
$ws.Range("C5").Value = "Ror2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 149.9875183333334
$ws.Range("H5").Value = 449.9625550000001
$ws.Range("I5").Value = 0.9582376829612175
$ws.Range("J5").Value = 0.9582376829612176
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1378523333333333
$ws.Range("N5").Value = 0.413557
$ws.Range("O5").Value = 0.01719056794796269
$ws.Range("P5").Value = 0.01719056794796269
$ws.Range("Q5").Value = 20.67612937312611
$ws.Range("R5").Value = 186.085164358135
$ws.Range("S5").Value = 0.01647264999924314
$ws.Range("T5").Value = 0.01647264999924315
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Cthrc1"
$ws.Range("C6").Value = "Ror2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 149.9875183333334
$ws.Range("H6").Value = 449.9625550000001
$ws.Range("I6").Value = 0.9582376829612175
$ws.Range("J6").Value = 0.9582376829612176
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 7.692787333333334
$ws.Range("N6").Value = 23.078362
$ws.Range("O6").Value = 0.9593118967607375
$ws.Range("P6").Value = 0.9593118967607375
$ws.Range("Q6").Value = 1153.822081192768
$ws.Range("R6").Value = 10384.39873073491
$ws.Range("S6").Value = 0.9192488091891398
$ws.Range("T6").Value = 0.9192488091891399
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Cthrc1"
$ws.Range("C7").Value = "Ror2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 149.9875183333334
$ws.Range("H7").Value = 449.9625550000001
$ws.Range("I7").Value = 0.9582376829612175
$ws.Range("J7").Value = 0.9582376829612176
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.1884283333333333
$ws.Range("N7").Value = 0.565285
$ws.Range("O7").Value = 0.02349753529129985
$ws.Range("P7").Value = 0.02349753529129985
$ws.Range("Q7").Value = 28.26189810035278
$ws.Range("R7").Value = 254.357082903175
$ws.Range("S7").Value = 0.02251622377283461
$ws.Range("T7").Value = 0.02251622377283461
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Cthrc1"
$ws.Range("C8").Value = "Ror2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.028937
$ws.Range("H8").Value = 3.086811
$ws.Range("I8").Value = 0.006573655046427582
$ws.Range("J8").Value = 0.006573655046427582
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1378523333333333
$ws.Range("N8").Value = 0.413557
$ws.Range("O8").Value = 0.01719056794796269
$ws.Range("P8").Value = 0.01719056794796269
$ws.Range("Q8").Value = 0.141841366303
$ws.Range("R8").Value = 1.276572296727
$ws.Range("S8").Value = 0.0001130048637420812
$ws.Range("T8").Value = 0.0001130048637420812
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Cthrc1"
$ws.Range("C9").Value = "Ror2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.028937
$ws.Range("H9").Value = 3.086811
$ws.Range("I9").Value = 0.006573655046427582
$ws.Range("J9").Value = 0.006573655046427582
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 7.692787333333334
$ws.Range("N9").Value = 23.078362
$ws.Range("O9").Value = 0.9593118967607375
$ws.Range("P9").Value = 0.9593118967607375
$ws.Range("Q9").Value = 7.915393520398001
$ws.Range("R9").Value = 71.238541683582
$ws.Range("S9").Value = 0.006306185491239237
$ws.Range("T9").Value = 0.006306185491239237
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Cthrc1"
$ws.Range("C10").Value = "Ror2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.028937
$ws.Range("H10").Value = 3.086811
$ws.Range("I10").Value = 0.006573655046427582
$ws.Range("J10").Value = 0.006573655046427582
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.1884283333333333
$ws.Range("N10").Value = 0.565285
$ws.Range("O10").Value = 0.02349753529129985
$ws.Range("P10").Value = 0.02349753529129985
$ws.Range("Q10").Value = 0.193880884015
$ws.Range("R10").Value = 1.744927956135
$ws.Range("S10").Value = 0.0001544646914462635
$ws.Range("T10").Value = 0.0001544646914462635
